# Add four new sheets after "Resultados Tests":
#   Resultados Tests1, Resultados Tests2, Resultados Tests3, NKAFFHA_valores
# Each "Resultados TestsN" sheet duplicates the small BET test-results table;
# NKAFFHA_valores combines the tail rows of the NKA and FFHA sheets.

$wb = $excel.ActiveWorkbook

function Fill-ResultadosTests($ws) {
    $ws.Range("A1").Value = "Test"
    $ws.Range("B1").Value = "Resultado"
    $ws.Range("C1").Value = "Promedio_A"
    $ws.Range("D1").Value = "Promedio_B"
    $ws.Range("E1").Value = "División"

    $ws.Range("A1:E1").Font.Bold = $true
    $ws.Range("A1:E1").HorizontalAlignment = -4108
    $ws.Range("A1:E1").VerticalAlignment = -4160
    $ws.Range("A1:E1").Borders.LineStyle = 1

    $ws.Range("A2").Value = "BET_BI"
    $ws.Range("B2").Value = "No hay poros cuello de botella"
    $ws.Range("C2").Value = "-"
    $ws.Range("D2").Value = "-"
    $ws.Range("E2").Value = "-"

    $ws.Range("A3").Value = "BET_P"
    $ws.Range("B3").Value = "Hay poros planos"
    $ws.Range("C3").Value = "-"
    $ws.Range("D3").Value = "-"
    $ws.Range("E3").Value = "-"

    $ws.Range("A4").Value = "BET_C"
    $ws.Range("B4").Value = "Hay poros cilindricos"
    $ws.Range("C4").Value = "-"
    $ws.Range("D4").Value = "-"
    $ws.Range("E4").Value = "-"
}

# --- Resultados Tests1 ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws1.Name = "Resultados Tests1"
Fill-ResultadosTests $ws1

# --- Resultados Tests2 ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws2.Name = "Resultados Tests2"
Fill-ResultadosTests $ws2

# --- Resultados Tests3 ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws3.Name = "Resultados Tests3"
Fill-ResultadosTests $ws3

# --- NKAFFHA_valores ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws4.Name = "NKAFFHA_valores"

$ws4.Range("A1").Value = "Radius of curvature"
$ws4.Range("B1").Value = "Vapor-Liquid Intrface Area"
$ws4.Range("C1").Value = "log(log(P/Po))"
$ws4.Range("D1").Value = "log(Vads)"

$ws4.Range("A1:D1").Font.Bold = $true
$ws4.Range("A1:D1").HorizontalAlignment = -4108
$ws4.Range("A1:D1").VerticalAlignment = -4160
$ws4.Range("A1:D1").Borders.LineStyle = 1

$ws4.Range("A2").Value = 4.1542
$ws4.Range("B2").Value = 9.578099999999999

$ws4.Range("A3").Value = 3.3475
$ws4.Range("B3").Value = 10.87

$ws4.Range("C4").Value = -0.0010658
$ws4.Range("D4").Value = 0.21327

$ws4.Range("C5").Value = 0.092696
$ws4.Range("D5").Value = 0.17002
